$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range('D2').Value = '61.605.84'
$ws.Range('E2').Value = '  -2.89%  '
$ws.Range('D3').Value = '2.577.62'
$ws.Range('E3').Value = '  -5.26%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '552.55'
$ws.Range('E5').Value = '  -1.39%  '
$ws.Range('D6').Value = '154.19'
$ws.Range('E6').Value = '  -2.42%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.596'
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('E9').Value = '  -3.00%  '
$ws.Range('D11').Value = '5.45'
$ws.Range('E11').Value = '  -3.65%  '
$ws.Range('D12').Value = '0.363'
$ws.Range('E12').Value = '  -2.80%  '
$ws.Range('D13').Value = '3.036.73'
$ws.Range('E13').Value = '  -5.08%  '
$ws.Range('D14').Value = '25.38'
$ws.Range('E14').Value = '  -4.32%  '
$ws.Range('D15').Value = '61.504.06'
$ws.Range('E15').Value = '  -2.84%  '
$ws.Range('E16').Value = '  -2.73%  '
$ws.Range('D17').Value = '2.581.10'
$ws.Range('E17').Value = '  -5.10%  '
$ws.Range('D18').Value = '11.55'
$ws.Range('E18').Value = '  -5.60%  '
$ws.Range('E19').Value = '  -2.99%  '
$ws.Range('D20').Value = '337.44'
$ws.Range('E20').Value = '  -3.88%  '
$ws.Range('E21').Value = '  -6.19%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '0.493'
$ws.Range('E23').Value = '  -4.18%  '
$ws.Range('D24').Value = '62.94'
$ws.Range('E24').Value = '  -2.28%  '
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').Value = '8.04'
$ws.Range('E27').Value = '  -2.14%  '
$ws.Range('D28').Value = '0.0₃0835'
$ws.Range('E28').Value = '  -5.38%  '
$ws.Range('E29').Value = '  -1.75%  '
$ws.Range('D30').Value = '7.04'
$ws.Range('E30').Value = '  -1.48%  '
$ws.Range('E31').Value = '  -4.89%  '
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').Value = '158.93'
$ws.Range('E33').Value = '  -4.22%  '
$ws.Range('D34').Value = '19.18'
$ws.Range('E34').Value = '  -3.49%  '
$ws.Range('D35').Value = '4.66'
$ws.Range('E35').Value = '  -3.40%  '
$ws.Range('E36').Value = '  -5.26%  '
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('D38').Value = '334.70'
$ws.Range('E38').Value = '  -3.15%  '
$ws.Range('D39').Value = '6.02'
$ws.Range('E39').Value = '  -0.79%  '
$ws.Range('D40').Value = '0.893'
$ws.Range('E40').Value = '  -7.04%  '
$ws.Range('D41').Value = '3.94'
$ws.Range('E41').Value = '  -2.90%  '
$ws.Range('D42').Value = '37.42'
$ws.Range('E42').Value = '  -2.08%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.145.27'
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '20.42'
$ws.Range('E45').Value = '  -4.56%  '
$ws.Range('D46').Value = '0.607'
$ws.Range('E46').Value = '  -3.05%  '
$ws.Range('D47').Value = '10.93'
$ws.Range('E47').Value = '  -1.24%  '
$ws.Range('E48').Value = '  -4.52%  '
$ws.Range('E49').Value = '  -6.12%  '
$ws.Range('D50').Value = '0.0966'
$ws.Range('E50').Value = '  -1.94%  '
$ws.Range('E51').Value = '  -2.63%  '
